$wb = $excel.ActiveWorkbook

# --- Sheet 1: Summary ---
$ws1 = $wb.Worksheets.Item("Summary")
$ws1.Range("B2").Value = 0.6226591760299626
$ws1.Range("C2").Value = 0.5771495877502945
$ws1.Range("D2").Value = 0.9176029962546817
$ws1.Range("E2").Value = 0.7086044830079538
$ws1.Range("F2").Value = 0.8207705192629816
$ws1.Range("G2").Value = 0.8972462849496443
$ws1.Range("H2").Value = 0.7754176661195977
$ws1.Range("I2").Value = 490
$ws1.Range("J2").Value = 359
$ws1.Range("K2").Value = 175
$ws1.Range("L2").Value = 44

# --- Sheet 2: Classification Report ---
$ws2 = $wb.Worksheets.Item("Classification Report")

# Row 2 ("0")
$ws2.Range("B2").Value = 0.7990867579908676
$ws2.Range("C2").Value = 0.3277153558052435
$ws2.Range("D2").Value = 0.4648074369189907

# Row 3 ("1")
$ws2.Range("B3").Value = 0.5771495877502945
$ws2.Range("C3").Value = 0.9176029962546817
$ws2.Range("D3").Value = 0.7086044830079538

# Row 4 ("accuracy")
$ws2.Range("B4").Value = 0.6226591760299626
$ws2.Range("C4").Value = 0.6226591760299626
$ws2.Range("D4").Value = 0.6226591760299626
$ws2.Range("E4").Value = 0.6226591760299626

# Row 5 ("macro avg")
$ws2.Range("B5").Value = 0.688118172870581
$ws2.Range("C5").Value = 0.6226591760299626
$ws2.Range("D5").Value = 0.5867059599634722

# Row 6 ("weighted avg")
$ws2.Range("B6").Value = 0.6881181728705811
$ws2.Range("C6").Value = 0.6226591760299626
$ws2.Range("D6").Value = 0.5867059599634722

# --- Sheet 3: Confusion Matrix ---
$ws3 = $wb.Worksheets.Item("Confusion Matrix")

# Row 2 ("Actual 0")
$ws3.Range("B2").Value = 175
$ws3.Range("C2").Value = 359

# Row 3 ("Actual 1")
$ws3.Range("B3").Value = 44
$ws3.Range("C3").Value = 490
